$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")

# Add new row 6 data to Sheet2
$ws2.Range("A6").Value = 6
$ws2.Range("B6").Value = "N"
$ws2.Range("C6").Value = "Wow, that's way off."

# Update selection on Sheet2 to C6 and make it the active/selected sheet
$ws2.Activate()
$ws2.Range("C6").Select()
